$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.107.38"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.740.75"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "3.739.20"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "4.366.11"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "3.737.77"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "69.043.14"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "3.886.07"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.108"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.673.28"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.93%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "423.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "2.779.85"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
